# The deck ships two DrawingML themes:
#   ppt/theme/theme1.xml -> "Office Theme" (default Office color scheme)
#   ppt/theme/theme2.xml -> "Integral" / "Red Violet" (the theme actually
#                            applied to the slide master / presentation)
#
# The authored change swaps which theme is "active": the presentation's
# applied color scheme becomes the default Office colors (what used to live
# in theme1.xml) instead of the Red Violet / Integral palette. Reproduce
# that by rewriting the active theme's 12-slot color scheme (reachable via
# ThemeColorScheme on the slide master's theme) to the stock "Office"
# values, in the standard ppThemeColorIndex order:
#   1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
#   5-10 Accent1-Accent6, 11 Hyperlink, 12 FollowedHyperlink

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# RRGGBB -> the BGR-packed long that PowerPoint's ColorFormat.RGB expects.
function HexToRgbLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @("000000","FFFFFF","44546A","E7E6E6","5B9BD5","ED7D31","A5A5A5","FFC000","4472C4","70AD47","0563C1","954F72")

for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = HexToRgbLong $officeColors[$i - 1]
}
